# consolidated_report.xlsx edit:
#   - add a new "giftcard_json_comparator" detector column (H), pushing the
#     existing "Final Result" (H->I) and "Reason" (I->J) columns one to the
#     right
#   - giftcard_consumer (C) now fails for every row
#   - recompute Final Result (I) / Reason (J) from the new json comparator
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at H. This shifts the old H ("Final Result") to I
# and the old I ("Reason") to J, and also picks up the bordered/bold header
# style for the new H1 automatically. Dimension grows from A1:I13 to A1:J13.
$ws.Columns("H").Insert()

# --- Header row ---
$ws.Range("H1").Value = "giftcard_json_comparator"
$ws.Range("I1").Value = "Final Result"
$ws.Range("J1").Value = "Reason"

# giftcard_consumer (C) fails across the board now
foreach ($r in 2..13) {
    $ws.Cells.Item($r, 3).Value = "Fail"
}

$jsonDivisionReason = "giftcard_json_comparator=[Reason=Field root['StyleHeaders'][0]['division'] exists in expected but missing in observed.; Field root['StyleHeaders'][0]['StyleDetails'][0]['division'] exists in expected but missing in observed.]"

# Per-row results: giftcard_json_comparator (H), Final Result (I), Reason (J)
$rowData = @{
    2  = @("Pass", "Fail", "")
    3  = @("Pass", "Fail", "")
    4  = @("Pass", "Fail", "")
    5  = @("Pass", "Fail", "")
    6  = @("Pass", "Fail", "")
    7  = @("Fail", "Fail", $jsonDivisionReason)
    8  = @("Fail", "Fail", $jsonDivisionReason)
    9  = @("Fail", "Fail", $jsonDivisionReason)
    10 = @("Fail", "Fail", $jsonDivisionReason)
    11 = @("Fail", "Fail", $jsonDivisionReason)
    12 = @("Pass", "Fail", "")
    13 = @("Pass", "Fail", "")
}

foreach ($r in 2..13) {
    $vals = $rowData[$r]
    $ws.Range("H$r").Value = $vals[0]
    $ws.Range("I$r").Value = $vals[1]
    $ws.Range("J$r").Value = $vals[2]
}
